$d = $word.ActiveDocument

# Paragraph 1: merge the split "This picture" / " " runs preceding the first
# image, and the split " " / "is an identicon." runs following it. A single
# text-preserving Find/Replace anywhere in the paragraph is enough to make
# Word re-coalesce all adjacent same-formatted runs in that paragraph.
$d.Content.Find.Execute("This picture ", $true, $false, $false, $false, $false, $true, 1, $false, "This picture ", 2)

# Paragraph 2: merge the split "Here is" / " " runs, the "one" / " " runs
# and " " / "that" runs inside the hyperlink, and the " " / "links." runs
# following the hyperlink. Same single-trigger-per-paragraph approach.
$d.Content.Find.Execute("Here is ", $true, $false, $false, $false, $false, $true, 1, $false, "Here is ", 2)
